# Correção de algumas funções e ajuste para rodar o fluxo novamente
#
# - Remove o 6o no (bus_006 / ip_006 / MAC004000*) da planilha "Load"
#   (linha 7 inteira) e limpa a linha correspondente em
#   "Public_Ilumination" (mantendo a linha, mas sem dados), reduzindo
#   o modelo de volta para 5 ns.
# - Corrige o horrio de start/end na aba "General" (B3).
# - Ajusta a aba/seleo ativa e as clulas selecionadas em cada planilha.

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General")
$wsBESS    = $wb.Worksheets.Item("BESS")
$wsGen     = $wb.Worksheets.Item("Generator")
$wsLoad    = $wb.Worksheets.Item("Load")
$wsPub     = $wb.Worksheets.Item("Public_Ilumination")

# --- General: corrige end_date (B3) de 06:30 para 06:00 ---
$wsGeneral.Range("B3").Value = 41098.25

# --- Load: remove os dados da linha 7 (6o n - bus_006 / MAC004000) ---
$wsLoad.Range("A7:L7").ClearContents()

# --- Public_Ilumination: limpa a linha 7 (6o n - ip_006 / bus_006) ---
$wsPub.Range("A7:L7").ClearContents()

# --- Seleo/aba ativa em cada planilha ---
$wsBESS.Range("G7").Select()
$wsGen.Range("B2").Select()
$wsLoad.Range("I14").Select()
$wsPub.Range("L12").Select()

# General volta a ser a aba ativa (tabSelected) com C6 selecionada
$wsGeneral.Activate()
$wsGeneral.Range("C6").Select()

# Melhor esforo: geometria da janela (nem sempre persistida pelo host)
$win = $excel.ActiveWindow
$win.Left = -38510
$win.Top = -5700
$win.Width = 38620
$win.Height = 21820

Write-Output "edit applied"
